# The diff adds an explicit <w:pageBreakBefore w:val="0"/> to:
#   1) the pPr of every paragraph in the document body
#   2) the pPr of the built-in Heading 1-6 / Title / Subtitle styles
# (the "Normal" / "Table Normal" styles are left untouched).
#
# Setting ParagraphFormat.PageBreakBefore explicitly (even to its current
# value of False/0) is what makes Word serialize the otherwise-implicit
# default as an explicit element, matching the target diff.

$d = $word.ActiveDocument

# 1) Every paragraph in the document body.
foreach ($p in $d.Paragraphs) {
    $p.Format.PageBreakBefore = 0
}

# 2) The relevant built-in paragraph styles.
$styleNames = @(
    "Heading 1",
    "Heading 2",
    "Heading 3",
    "Heading 4",
    "Heading 5",
    "Heading 6",
    "Title",
    "Subtitle"
)

foreach ($styleName in $styleNames) {
    $style = $d.Styles($styleName)
    $style.ParagraphFormat.PageBreakBefore = 0
}
